$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data row for "the north" / "miền bắc" (row 6: A6="the north", B6="miền bắc")
$ws.Rows.Item(6).Delete()
